$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D="67.683.30"; DText=$false; E="  +0.35%  "},
    @{Row=3; D="3.301.69"; DText=$false; E="  -2.02%  "},
    @{Row=4; D="1.00"; DText=$true; E="  +0.13%  "},
    @{Row=5; D="581.38"; DText=$true; E="  -1.48%  "},
    @{Row=6; D="174.93"; DText=$true; E="  -6.96%  "},
    @{Row=7; D="1.00"; DText=$true; E="  +0.01%  "},
    @{Row=8; D="0.580"; DText=$true; E="  -3.00%  "},
    @{Row=9; D="3.298.70"; DText=$false; E="  -1.97%  "},
    @{Row=10; D="0.174"; DText=$true; E="  -5.23%  "},
    @{Row=11; D="0.572"; DText=$true; E="  -2.48%  "},
    @{Row=12; D="45.19"; DText=$true; E="  -4.77%  "},
    @{Row=13; E="  -2.47%  "},
    @{Row=14; D="666.42"; DText=$true; E="  +4.07%  "},
    @{Row=15; D="3.826.06"; DText=$false; E="  -2.30%  "},
    @{Row=16; D="8.33"; DText=$true; E="  -3.36%  "},
    @{Row=17; D="67.764.59"; DText=$false; E="  +0.41%  "},
    @{Row=18; E="  -0.46%  "},
    @{Row=19; D="3.297.33"; DText=$false; E="  -2.28%  "},
    @{Row=20; D="17.38"; DText=$true; E="  -3.60%  "},
    @{Row=21; D="10.83"; DText=$true; E="  -3.26%  "},
    @{Row=22; D="0.884"; DText=$true; E="  -2.95%  "},
    @{Row=23; D="5.39"; DText=$true; E="  +5.87%  "},
    @{Row=24; D="16.99"; DText=$true; E="  -5.79%  "},
    @{Row=25; D="98.43"; DText=$true; E="  -1.89%  "},
    @{Row=26; E="  -3.94%  "},
    @{Row=27; D="2.65"; DText=$true; E="  -6.71%  "},
    @{Row=28; D="9.16"; DText=$true; E="  -5.89%  "},
    @{Row=29; D="32.82"; DText=$true; E="  +1.06%  "},
    @{Row=30; D="8.34"; DText=$true; E="  -4.21%  "},
    @{Row=31; D="6.98"; DText=$true; E="  +1.42%  "},
    @{Row=32; D="575.17"; DText=$true; E="  -6.08%  "},
    @{Row=33; D="10.91"; DText=$true; E="  -1.67%  "},
    @{Row=34; D="3.761.76"; DText=$false; E="  -4.08%  "},
    @{Row=35; D="0.103"; DText=$true; E="  -3.32%  "},
    @{Row=36; D="1.00"; DText=$true; E="  +0.12%  "},
    @{Row=37; D="3.36"; DText=$true; E="  -13.32%  "},
    @{Row=38; D="55.65"; DText=$true; E="  -0.66%  "},
    @{Row=39; E="  -1.81%  "},
    @{Row=40; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="2.62"; DText=$true; E="  -7.29%  "},
    @{Row=41; B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="32.28"; DText=$true; E="  -4.38%  "},
    @{Row=42; D="3.06"; DText=$true; E="  -7.51%  "},
    @{Row=43; D="0.0₃0660"; DText=$false; E="  -6.33%  "},
    @{Row=44; B="TheGraph"; C="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; D="0.327"; DText=$true; E="  -5.05%  "},
    @{Row=45; B="ApeXProtocol"; C="https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D="3.24"; DText=$true; E="  -5.07%  "},
    @{Row=46; D="0.0402"; DText=$true; E="  -5.22%  "},
    @{Row=47; E="  -0.44%  "},
    @{Row=48; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.127"; DText=$true; E="  -2.29%  "},
    @{Row=49; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="1.00"; DText=$true; E="  +0.02%  "},
    @{Row=50; E="  -0.86%  "},
    @{Row=51; D="2.76"; DText=$true; E="  -1.89%  "}
)

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) {
        $ws.Cells.Item($r.Row, 2).Value = $r.B
    }
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        if ($r.DText) {
            $ws.Cells.Item($r.Row, 4).NumberFormat = "@"
        }
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    if ($r.ContainsKey("E")) {
        $ws.Cells.Item($r.Row, 5).Value = $r.E
    }
}
